$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.996.87"
$ws.Range("E2").Value = "  +4.47%  "
$ws.Range("D3").Value = "2.467.25"
$ws.Range("E3").Value = "  +5.73%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'565.95"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "'143.20"
$ws.Range("E6").Value = "  +9.32%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("D9").Value = "2.466.27"
$ws.Range("E9").Value = "  +5.85%  "
$ws.Range("E10").Value = "  +3.74%  "
$ws.Range("D11").Value = "'5.69"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("E13").Value = "  +4.50%  "
$ws.Range("D14").Value = "'26.41"
$ws.Range("E14").Value = "  +11.44%  "
$ws.Range("D15").Value = "2.908.65"
$ws.Range("E15").Value = "  +5.88%  "
$ws.Range("D16").Value = "62.921.28"
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("E17").Value = "  +4.59%  "
$ws.Range("D18").Value = "2.468.06"
$ws.Range("E18").Value = "  +5.79%  "
$ws.Range("D19").Value = "'11.23"
$ws.Range("E19").Value = "  +5.01%  "
$ws.Range("D20").Value = "'341.03"
$ws.Range("E20").Value = "  +8.21%  "
$ws.Range("D21").Value = "'4.27"
$ws.Range("E21").Value = "  +4.27%  "
$ws.Range("E22").Value = "  +3.06%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'65.51"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("E27").Value = "  +8.70%  "
$ws.Range("D28").Value = "'8.06"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("E29").Value = "  +6.29%  "
$ws.Range("D30").Value = "'6.82"
$ws.Range("E30").Value = "  +12.31%  "
$ws.Range("D31").Value = "0.0₃0803"
$ws.Range("E31").Value = "  +9.56%  "
$ws.Range("E32").Value = "  +6.41%  "
$ws.Range("D33").Value = "'177.68"
$ws.Range("E33").Value = "  +3.58%  "
$ws.Range("E34").Value = "  +11.02%  "
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("D36").Value = "'18.87"
$ws.Range("E36").Value = "  +4.36%  "
$ws.Range("D37").Value = "'367.23"
$ws.Range("E37").Value = "  +13.77%  "
$ws.Range("D38").Value = "'4.41"
$ws.Range("E38").Value = "  +7.32%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +10.86%  "
$ws.Range("D42").Value = "'40.48"
$ws.Range("E42").Value = "  +6.20%  "
$ws.Range("D43").Value = "'150.48"
$ws.Range("E43").Value = "  +9.19%  "
$ws.Range("D44").Value = "'3.70"
$ws.Range("E44").Value = "  +5.43%  "
$ws.Range("D45").Value = "'20.54"
$ws.Range("E45").Value = "  +6.04%  "
$ws.Range("D46").Value = "'0.597"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("D47").Value = "'0.0958"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").Value = "'0.0516"
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("D49").Value = "0.0₆0239"
$ws.Range("E49").Value = "  +9.25%  "
$ws.Range("E50").Value = "  +4.28%  "
$ws.Range("D51").Value = "'17.98"
$ws.Range("E51").Value = "  +5.36%  "
